$d = $word.ActiveDocument

# --- 1) Title ---
$d.Content.Find.Execute("Transcendence of Artistic Shadows", $true, $false, $false, $false, $false, $true, 1, $false, "Exploring the Human Body: An Introduction to Biology for High School Students", 2) | Out-Null

# --- 2) Author name ---
$d.Content.Find.Execute("Gilbert James", $true, $false, $false, $false, $false, $true, 1, $false, "Cassandra Lee", 2) | Out-Null

# --- 3) Email username (domain ".com" stays the same via a separate run) ---
$d.Content.Find.Execute("gilbertjames@figaromail", $true, $false, $false, $false, $false, $true, 1, $false, "cassandralee@validmail", 2) | Out-Null

# --- 4) Body intro paragraph (paragraph 5) -> full rewrite ---
$p5 = $d.Paragraphs.Item(5)
$p5xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" />
          <w:color w:val="000000" />
          <w:sz w:val="24" />
        </w:rPr>
        <w:t>Biology, the study of life, is a vast and complex discipline that encompasses everything from the tiniest molecules to the largest organisms</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" />
          <w:color w:val="000000" />
          <w:sz w:val="24" />
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" />
          <w:color w:val="000000" />
          <w:sz w:val="24" />
        </w:rPr>
        <w:t xml:space="preserve"> In this essay, we will embark on a journey to explore the wonders of human biology, learning about the systems that make us function and the incredible processes that allow us to thrive</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" />
          <w:color w:val="000000" />
          <w:sz w:val="24" />
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" />
          <w:color w:val="000000" />
          <w:sz w:val="24" />
        </w:rPr>
        <w:br />
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" />
          <w:color w:val="000000" />
          <w:sz w:val="24" />
        </w:rPr>
        <w:br />
        <w:t>Paragraph 1:</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" />
          <w:color w:val="000000" />
          <w:sz w:val="24" />
        </w:rPr>
        <w:br />
        <w:t>Life is a remarkable phenomenon, and human beings are among the most intricate and fascinating organisms in the natural world</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" />
          <w:color w:val="000000" />
          <w:sz w:val="24" />
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" />
          <w:color w:val="000000" />
          <w:sz w:val="24" />
        </w:rPr>
        <w:t xml:space="preserve"> Comprised of trillions of cells, each with its own unique function, our bodies are capable of astounding feats of adaptation, growth, and reproduction</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" />
          <w:color w:val="000000" />
          <w:sz w:val="24" />
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" />
          <w:color w:val="000000" />
          <w:sz w:val="24" />
        </w:rPr>
        <w:t xml:space="preserve"> From the intricate workings of our cells to the complex interactions between our organs, human biology is a story of interconnectedness and resilience</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" />
          <w:color w:val="000000" />
          <w:sz w:val="24" />
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" />
          <w:color w:val="000000" />
          <w:sz w:val="24" />
        </w:rPr>
        <w:br />
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" />
          <w:color w:val="000000" />
          <w:sz w:val="24" />
        </w:rPr>
        <w:br />
        <w:t>Paragraph 2:</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" />
          <w:color w:val="000000" />
          <w:sz w:val="24" />
        </w:rPr>
        <w:br />
        <w:t>Biologists seek to understand the fundamental principles that govern life, examining everything from the structure of molecules to the behavior of entire ecosystems</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" />
          <w:color w:val="000000" />
          <w:sz w:val="24" />
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" />
          <w:color w:val="000000" />
          <w:sz w:val="24" />
        </w:rPr>
        <w:t xml:space="preserve"> Whether studying the genetic code that shapes our traits or investigating the intricate web of interactions that occur within cells, biologists are constantly striving to unravel the mysteries of the living world</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" />
          <w:color w:val="000000" />
          <w:sz w:val="24" />
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" />
          <w:color w:val="000000" />
          <w:sz w:val="24" />
        </w:rPr>
        <w:t xml:space="preserve"> Their discoveries have led to groundbreaking advancements in medicine, agriculture, and environmental science, benefitting humanity in countless ways</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" />
          <w:color w:val="000000" />
          <w:sz w:val="24" />
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" />
          <w:color w:val="000000" />
          <w:sz w:val="24" />
        </w:rPr>
        <w:br />
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" />
          <w:color w:val="000000" />
          <w:sz w:val="24" />
        </w:rPr>
        <w:br />
        <w:t>Paragraph 3:</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" />
          <w:color w:val="000000" />
          <w:sz w:val="24" />
        </w:rPr>
        <w:br />
        <w:t>As high school students, we have a unique opportunity to delve into the fascinating world of biology</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" />
          <w:color w:val="000000" />
          <w:sz w:val="24" />
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" />
          <w:color w:val="000000" />
          <w:sz w:val="24" />
        </w:rPr>
        <w:t xml:space="preserve"> Through classroom lessons, field trips, and hands-on experiments, we can gain a deeper understanding of the human body and its amazing capabilities</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" />
          <w:color w:val="000000" />
          <w:sz w:val="24" />
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" />
          <w:color w:val="000000" />
          <w:sz w:val="24" />
        </w:rPr>
        <w:t xml:space="preserve"> Exploring topics like cell division, genetics, and anatomy can help us appreciate the complexity of life and the interconnectedness of all living things</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" />
          <w:color w:val="000000" />
          <w:sz w:val="24" />
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" />
          <w:color w:val="000000" />
          <w:sz w:val="24" />
        </w:rPr>
        <w:t xml:space="preserve"> Moreover, biology offers an avenue for developing critical thinking and problem-solving skills that will serve us well in any future endeavor</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" />
          <w:color w:val="000000" />
          <w:sz w:val="24" />
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
    </w:p>
'@
$p5.Range.InsertXML($p5xml)

# --- 5) "Summary" heading paragraph -> add lastRenderedPageBreak ---
$p6 = $d.Paragraphs.Item(6)
$p6xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" />
          <w:color w:val="000000" />
          <w:sz w:val="28" />
        </w:rPr>
        <w:lastRenderedPageBreak />
        <w:t>Summary</w:t>
      </w:r>
    </w:p>
'@
$p6.Range.InsertXML($p6xml)

# --- 6) Conclusion paragraph -> full rewrite ---
$p7 = $d.Paragraphs.Item(7)
$p7xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" />
          <w:color w:val="000000" />
        </w:rPr>
        <w:t>In this essay, we provided a brief introduction to biology, highlighting its significance in understanding human life and the natural world</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" />
          <w:color w:val="000000" />
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" />
          <w:color w:val="000000" />
        </w:rPr>
        <w:t xml:space="preserve"> We discussed the interconnectedness of living organisms and the astounding complexity of the human body</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" />
          <w:color w:val="000000" />
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" />
          <w:color w:val="000000" />
        </w:rPr>
        <w:t xml:space="preserve"> Moreover, we emphasized the importance of biological research in advancing medicine and agriculture</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" />
          <w:color w:val="000000" />
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" />
          <w:color w:val="000000" />
        </w:rPr>
        <w:t xml:space="preserve"> By exploring the wonders of biology, high school students can develop a deeper appreciation for life and acquire essential skills that will benefit them in their future pursuits</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" />
          <w:color w:val="000000" />
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
    </w:p>
'@
$p7.Range.InsertXML($p7xml)

# --- 7) Append a new empty paragraph at the end of the document ---
$d.Content.InsertParagraphAfter()

Write-Output "done"
